$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 4405
$ws.Range("D2").Value = 489
$ws.Range("E2").Value = 1272
$ws.Range("H2").Value = 2137
$ws.Range("I2").Value = 16510
$ws.Range("J2").Value = 7
$ws.Range("K2").Value = 693
$ws.Range("L2").Value = 20014
$ws.Range("D3").Value = 1.237
$ws.Range("E3").Value = 1.242
$ws.Range("H3").Value = 8
$ws.Range("I3").Value = 7461
$ws.Range("J3").Value = 22
$ws.Range("B4").Value = 263
$ws.Range("B5").Value = 804
$ws.Range("D5").Value = 407
$ws.Range("E5").Value = 411
$ws.Range("H5").Value = 73
$ws.Range("I5").Value = 7739
$ws.Range("B6").Value = 264
$ws.Range("B7").Value = 260
$ws.Range("D7").Value = 19
$ws.Range("E7").Value = 19
$ws.Range("I7").Value = 6316
$ws.Range("B8").Value = 776
$ws.Range("D8").Value = 265
$ws.Range("E8").Value = 275
$ws.Range("H8").Value = 257
$ws.Range("I8").Value = 6679
$ws.Range("B9").Value = 778
$ws.Range("D9").Value = 82
$ws.Range("E9").Value = 83
$ws.Range("H9").Value = 120
$ws.Range("I9").Value = 9268
$ws.Range("B10").Value = 756
$ws.Range("D10").Value = 83
$ws.Range("E10").Value = 85
$ws.Range("H10").Value = 217
$ws.Range("I10").Value = 12339
$ws.Range("B11").Value = 928
$ws.Range("D11").Value = 251
$ws.Range("E11").Value = 258
$ws.Range("H11").Value = 2825
$ws.Range("I11").Value = 41975
$ws.Range("J11").Value = 7
$ws.Range("B12").Value = 787
$ws.Range("D12").Value = 39
$ws.Range("E12").Value = 38
$ws.Range("H12").Value = 313
$ws.Range("I12").Value = 14970
$ws.Range("B13").Value = 561
$ws.Range("D13").Value = 245
$ws.Range("E13").Value = 279
$ws.Range("F13").Value = 37
$ws.Range("H13").Value = 5191
$ws.Range("I13").Value = 9266
$ws.Range("B14").Value = 805
$ws.Range("D14").Value = 308
$ws.Range("E14").Value = 310
$ws.Range("H14").Value = 32
$ws.Range("I14").Value = 7815
$ws.Range("J14").Value = 7
$ws.Range("B15").Value = 262
$ws.Range("B16").Value = 303
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 439
$ws.Range("E16").Value = 633
$ws.Range("H16").Value = 5391
$ws.Range("I16").Value = 8691
$ws.Range("K16").Value = 118
$ws.Range("L16").Value = 2538
$ws.Range("B17").Value = 51
$ws.Range("D17").Value = 30
$ws.Range("E17").Value = 30
$ws.Range("I17").Value = 5000
$ws.Range("B18").Value = 221
$ws.Range("D18").Value = 596
$ws.Range("E18").Value = 716
$ws.Range("H18").Value = 591
$ws.Range("I18").Value = 5464
$ws.Range("J18").Value = 11
$ws.Range("K18").Value = 74
$ws.Range("L18").Value = 1034
$ws.Range("B19").Value = 766
$ws.Range("D19").Value = 237
$ws.Range("E19").Value = 267
$ws.Range("H19").Value = 2265
$ws.Range("I19").Value = 7832
$ws.Range("B20").Value = 1394
$ws.Range("B21").Value = 437
$ws.Range("D21").Value = 148
$ws.Range("E21").Value = 156
$ws.Range("H21").Value = 1584
$ws.Range("I21").Value = 21196
$ws.Range("J21").Value = 5
$ws.Range("B22").Value = 183
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = 8
$ws.Range("H22").Value = 3750
$ws.Range("I22").Value = 2000
$ws.Range("B23").Value = 814
$ws.Range("D23").Value = 21
$ws.Range("E23").Value = 21
$ws.Range("I23").Value = 8095
$ws.Range("B24").Value = 957
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 15
$ws.Range("H24").Value = 1250
$ws.Range("I24").Value = 7143
$ws.Range("J24").Value = 1
$ws.Range("B25").Value = 63
$ws.Range("D25").Value = 76
$ws.Range("E25").Value = 101
$ws.Range("H25").Value = 396
$ws.Range("I25").Value = 3553
$ws.Range("L25").Value = 2079
$ws.Range("B26").Value = 123
$ws.Range("D26").Value = 51
$ws.Range("E26").Value = 51
$ws.Range("I26").Value = 8824
$ws.Range("B27").Value = 179
$ws.Range("D27").Value = 58
$ws.Range("E27").Value = 75
$ws.Range("H27").Value = 2267
$ws.Range("I27").Value = 3103
